# 09.12.2020 MC Sales Details
# Swap the Old/New retail name values on the existing row (RET-36547),
# add a new retailer name-change entry (RET-36436: Ma Telecom -> Jewel Telecom),
# and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap B2 (Old Retail Name) and C2 (New Retail Name) values for RET-36547
$oldB2 = $ws.Range("B2").Value2
$oldC2 = $ws.Range("C2").Value2
$ws.Range("B2").Value = $oldC2
$ws.Range("C2").Value = $oldB2

# Add new retailer name-change row for RET-36436
$ws.Range("A3").Value = "RET-36436"
$ws.Range("B3").Value = "Ma Telecom"
$ws.Range("C3").Value = "Jewel Telecom"
$ws.Range("D3").Value = "N/A"
$ws.Range("E3").Value = "N/A"
$ws.Range("F3").Value = "N/A"

# Update the active cell selection
$ws.Range("C13").Select()
